$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5000.0835
$ws.Range("J40").Value = 5150.1
$ws.Range("L40").Value = 5150.1
$ws.Range("N40").Value = -5500.1
# Row 62
$ws.Range("H62").Value = 37601.7
$ws.Range("I62").Value = 39833
$ws.Range("J62").Value = 36645.43
$ws.Range("K62").Value = 39833
$ws.Range("L62").Value = 36645.43
$ws.Range("M62").Value = -39209
$ws.Range("N62").Value = -37893.43
# Row 65
$ws.Range("H65").Value = 37601.7
$ws.Range("I65").Value = 39833
$ws.Range("J65").Value = 36645.43
$ws.Range("K65").Value = 199165
$ws.Range("L65").Value = 183227.15
$ws.Range("M65").Value = -196045
$ws.Range("N65").Value = -189467.15
# Row 88
$ws.Range("H88").Value = 113212.664
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 113212.664
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 113212.664
$ws.Range("N88").Value = -114024.664
# Row 91
$ws.Range("H91").Value = 113212.664
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 113212.664
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 113212.664
$ws.Range("N91").Value = -116020.664
# Row 98
$ws.Range("H98").Value = 1981.8148
$ws.Range("I98").Value = 1771.8334
$ws.Range("J98").Value = 2401.7778
$ws.Range("K98").Value = 1771.8334
$ws.Range("L98").Value = 2401.7778
$ws.Range("M98").Value = -273.8334
$ws.Range("N98").Value = -5397.7778
# Row 101
$ws.Range("H101").Value = 3333910.8
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 106
$ws.Range("H106").Value = 1304.8462
$ws.Range("I106").Value = 951.63635
$ws.Range("K106").Value = 951.63635
$ws.Range("M106").Value = -320.63635
# Row 113
$ws.Range("H113").Value = 6452.8
$ws.Range("I113").Value = 6212.25
$ws.Range("J113").Value = 6540.273
$ws.Range("K113").Value = 6212.25
$ws.Range("L113").Value = 6540.273
$ws.Range("M113").Value = -2958.25
$ws.Range("N113").Value = -13048.273
# Row 122
$ws.Range("H122").Value = 1981.8148
$ws.Range("I122").Value = 1771.8334
$ws.Range("J122").Value = 2401.7778
$ws.Range("K122").Value = 5315.5002
$ws.Range("L122").Value = 7205.3334
$ws.Range("M122").Value = -2865.5002
$ws.Range("N122").Value = -12105.3334
# Row 138
$ws.Range("H138").Value = 1326
$ws.Range("I138").Value = 1326
$ws.Range("K138").Value = 3978
$ws.Range("M138").Value = 1162

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3880270.8
$ws.Range("I32").Value = 746157.7
$ws.Range("K32").Value = 746157.7
$ws.Range("M32").Value = -745870.7
# Row 45
$ws.Range("H45").Value = 14051.56
$ws.Range("I45").Value = 12072.105
$ws.Range("J45").Value = 20319.834
$ws.Range("K45").Value = 12072.105
$ws.Range("L45").Value = 20319.834
$ws.Range("M45").Value = -11695.105
$ws.Range("N45").Value = -21073.834

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 774.75
$ws.Range("I22").Value = 899.6667
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 899.6667
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -726.6667
$ws.Range("N22").Value = -746
# Row 102
$ws.Range("H102").Value = 15293.714
$ws.Range("I102").Value = 12740.667
$ws.Range("K102").Value = 12740.667
$ws.Range("M102").Value = -9495.666999999999
# Row 107
$ws.Range("H107").Value = 2192.625
$ws.Range("J107").Value = 1644.6364
$ws.Range("L107").Value = 1644.6364
$ws.Range("N107").Value = -5484.6364

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 125000230
$ws.Range("I7").Value = 166666880
$ws.Range("J7").Value = 260
$ws.Range("K7").Value = 166666880
$ws.Range("L7").Value = 260
$ws.Range("M7").Value = -166666767
$ws.Range("N7").Value = -486
# Row 22
$ws.Range("H22").Value = 1234.1428
$ws.Range("I22").Value = 1197.8
$ws.Range("J22").Value = 1325
$ws.Range("K22").Value = 1197.8
$ws.Range("L22").Value = 1325
$ws.Range("M22").Value = -847.8
$ws.Range("N22").Value = -2025
# Row 60
$ws.Range("H60").Value = 9999
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 62
$ws.Range("H62").Value = 999
$ws.Range("I62").Value = 999
$ws.Range("K62").Value = 999
$ws.Range("M62").Value = -375
# Row 65
$ws.Range("H65").Value = 999
$ws.Range("I65").Value = 999
$ws.Range("K65").Value = 4995
$ws.Range("M65").Value = -1875
# Row 105
$ws.Range("H105").Value = 639.8
$ws.Range("I105").Value = 599.75
$ws.Range("K105").Value = 599.75
$ws.Range("M105").Value = 1147.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 49.068966
$ws.Range("I2").Value = 34.666668
$ws.Range("J2").Value = 72.63636
$ws.Range("K2").Value = 208.000008
$ws.Range("L2").Value = 435.81816
$ws.Range("M2").Value = -95.00000800000001
$ws.Range("N2").Value = -661.81816
# Row 110
$ws.Range("H110").Value = 12796
$ws.Range("J110").Value = 19996.666
$ws.Range("L110").Value = 59989.99800000001
$ws.Range("N110").Value = -68169.99800000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 42627.75
$ws.Range("J39").Value = 42627.75
$ws.Range("L39").Value = 42627.75
$ws.Range("N39").Value = -43691.75
# Row 97
$ws.Range("H97").Value = 65840.95
$ws.Range("J97").Value = 164310
$ws.Range("L97").Value = 164310
$ws.Range("N97").Value = -165302
# Row 122
$ws.Range("H122").Value = 2769.25
$ws.Range("I122").Value = 1285.875
$ws.Range("K122").Value = 3857.625
$ws.Range("M122").Value = -1407.625
# Row 132
$ws.Range("H132").Value = 4669.4
$ws.Range("J132").Value = 4071.2222
$ws.Range("L132").Value = 12213.6666
$ws.Range("N132").Value = -17273.6666

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1297.4
$ws.Range("I16").Value = 1347.1111
$ws.Range("K16").Value = 1347.1111
$ws.Range("M16").Value = -1177.1111
# Row 22
$ws.Range("H22").Value = 1232.591
$ws.Range("I22").Value = 788.75
$ws.Range("J22").Value = 1486.2142
$ws.Range("K22").Value = 788.75
$ws.Range("L22").Value = 1486.2142
$ws.Range("M22").Value = -493.75
$ws.Range("N22").Value = -2076.2142
# Row 27
$ws.Range("H27").Value = 1232.591
$ws.Range("I27").Value = 788.75
$ws.Range("J27").Value = 1486.2142
$ws.Range("K27").Value = 788.75
$ws.Range("L27").Value = 1486.2142
$ws.Range("M27").Value = -681.75
$ws.Range("N27").Value = -1700.2142
# Row 40
$ws.Range("H40").Value = 6212.0835
$ws.Range("I40").Value = 6050
$ws.Range("K40").Value = 6050
$ws.Range("M40").Value = -5914
# Row 55
$ws.Range("H55").Value = 338.09525
$ws.Range("I55").Value = 243.3125
$ws.Range("K55").Value = 243.3125
$ws.Range("M55").Value = -70.3125
# Row 61
$ws.Range("H61").Value = 88662.94
$ws.Range("I61").Value = 82187.84
$ws.Range("K61").Value = 82187.84
$ws.Range("M61").Value = -81985.84
# Row 82
$ws.Range("H82").Value = 4157.3335
$ws.Range("I82").Value = 4269.9
$ws.Range("K82").Value = 4269.9
$ws.Range("M82").Value = -3908.9
# Row 85
$ws.Range("H85").Value = 4157.3335
$ws.Range("I85").Value = 4269.9
$ws.Range("K85").Value = 4269.9
$ws.Range("M85").Value = -3021.9
# Row 113
$ws.Range("H113").Value = 88662.94
$ws.Range("I113").Value = 82187.84
$ws.Range("K113").Value = 82187.84
$ws.Range("M113").Value = -80017.84
# Row 122
$ws.Range("H122").Value = 4167.4287
$ws.Range("I122").Value = 3486
$ws.Range("K122").Value = 10458
$ws.Range("M122").Value = -8008

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 107
$ws.Range("H107").Value = 29446882
$ws.Range("I107").Value = 1881.0714
$ws.Range("K107").Value = 5643.2142
$ws.Range("M107").Value = -3723.2142
# Row 122
$ws.Range("H122").Value = 1525.5745
$ws.Range("I122").Value = 1114.9706
$ws.Range("K122").Value = 3344.9118
$ws.Range("M122").Value = -894.9118000000003
# Row 132
$ws.Range("H132").Value = 1202204.1
$ws.Range("I132").Value = 1468440.4
$ws.Range("K132").Value = 4405321.199999999
$ws.Range("M132").Value = -4402791.199999999

$wb.Save()
Write-Output "Applied Leviathan_Profits updates"